$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Cells changing from numeric style to shared-text style ("0" / "***.*") ---
# Use format-copy (style only) then value-copy (text value) from donor cells that
# already hold the desired text + style, so we reuse existing styles/shared strings
# instead of creating new ones.
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4163) | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null

$ws.Range("E14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").PasteSpecial(-4163) | Out-Null

# --- Cells changing from shared-text style to numeric style ---
$ws.Range("C18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("K18").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -88
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 500
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 80
$ws.Range("I15").Value = 18
$ws.Range("J15").Value = 21
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = 63.636363636363
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = -25
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 174
$ws.Range("J16").Value = 191
$ws.Range("K16").Value = -8.900523560209
$ws.Range("L16").Value = 74
$ws.Range("M16").Value = 18.367346938775
$ws.Range("N16").Value = -73.636363636363
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 70
$ws.Range("F17").Value = 70
$ws.Range("G17").Value = 57
$ws.Range("H17").Value = 22.807017543859
$ws.Range("I17").Value = 261
$ws.Range("J17").Value = 214
$ws.Range("K17").Value = 21.962616822429
$ws.Range("L17").Value = 42.622950819672
$ws.Range("M17").Value = 61.111111111111
$ws.Range("N17").Value = -1.87969924812
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 56.25
$ws.Range("I18").Value = 124
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = 36.263736263736
$ws.Range("L18").Value = 56.962025316455
$ws.Range("M18").Value = -0.8
$ws.Range("N18").Value = -78.359511343804
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 10.526315789473
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 83
$ws.Range("H19").Value = 3.614457831325
$ws.Range("I19").Value = 322
$ws.Range("J19").Value = 342
$ws.Range("K19").Value = -5.847953216374
$ws.Range("L19").Value = 64.285714285714
$ws.Range("M19").Value = 93.975903614457
$ws.Range("N19").Value = 29.317269076305
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -18.181818181818
$ws.Range("F20").Value = 57
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = 54.054054054054
$ws.Range("I20").Value = 216
$ws.Range("J20").Value = 190
$ws.Range("K20").Value = 13.684210526315
$ws.Range("L20").Value = 113.861386138614
$ws.Range("M20").Value = 213.04347826087
$ws.Range("N20").Value = -66.091051805337
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 51
$ws.Range("E21").Value = 27.450980392156
$ws.Range("F21").Value = 280
$ws.Range("G21").Value = 234
$ws.Range("H21").Value = 19.658119658119
$ws.Range("I21").Value = 1118
$ws.Range("J21").Value = 1053
$ws.Range("K21").Value = 6.172839506172
$ws.Range("L21").Value = 65.140324963072
$ws.Range("M21").Value = 64.170337738619
$ws.Range("N21").Value = -54.067378800328
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 26
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 90
$ws.Range("J23").Value = 108
$ws.Range("K23").Value = -16.666666666666
$ws.Range("L23").Value = 18.421052631578
$ws.Range("M23").Value = 32.35294117647
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = -25.581395348837
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 173
$ws.Range("H24").Value = -30.057803468208
$ws.Range("I24").Value = 640
$ws.Range("J24").Value = 609
$ws.Range("K24").Value = 5.090311986863
$ws.Range("L24").Value = 74.863387978142
$ws.Range("M24").Value = 40.969162995594
$ws.Range("C25").Value = 23
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 43.75
$ws.Range("F25").Value = 86
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 347
$ws.Range("J25").Value = 353
$ws.Range("K25").Value = -1.699716713881
$ws.Range("L25").Value = 29.962546816479
$ws.Range("M25").Value = -28.600823045267
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 600
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 120
$ws.Range("I26").Value = 24
$ws.Range("J26").Value = 30
$ws.Range("K26").Value = -20
$ws.Range("L26").Value = 41.176470588235
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 36
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = 71.428571428571
$ws.Range("L27").Value = 16.129032258064
$ws.Range("N28").Value = -77.777777777777
$ws.Range("N29").Value = -83.333333333333

